$d = $word.ActiveDocument

# The {{date}} placeholder paragraph is expanded into the full set of
# letter placeholders, each on its own paragraph, in order:
#   {{refNo}}, {{date}}, {{recipientName}}, {{recipientAddress}},
#   {{subject}}, {{content}}, {{senderName}}, {{senderPosition}},
#   {{organization}}
# A straightforward (dumb) find/replace does the job: replace the
# "{{date}}" text with the same token plus the rest of the tokens,
# each separated by a paragraph mark ("`r").
$replacement = "{{refNo}}`r{{date}}`r{{recipientName}}`r{{recipientAddress}}`r{{subject}}`r{{content}}`r{{senderName}}`r{{senderPosition}}`r{{organization}}"

$d.Content.Find.Execute("{{date}}", $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)
